$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 192 ---
$ws.Cells.Item(192, 1).Value = 45470.2916666667
$ws.Cells.Item(192, 2).Value = 0
$ws.Cells.Item(192, 3).Value = 4.1100001335144
$ws.Cells.Item(192, 4).Value = 4.1100001335144
$ws.Cells.Item(192, 5).Value = 4.1100001335144
$ws.Cells.Item(192, 6).Value = 4.1100001335144
# adj_close (col G) is stored as text in the source data even though it
# looks numeric - force text formatting before the write so it lands as a
# shared string rather than being coerced back to a number.
$ws.Cells.Item(192, 7).NumberFormat = "@"
$ws.Cells.Item(192, 7).Value = "4.1100001335144"
$ws.Cells.Item(192, 7).Style = "Normal"
$ws.Cells.Item(192, 8).Value = "ESF.MI"

# --- Row 193 ---
$ws.Cells.Item(193, 1).Value = 45471.5052199074
$ws.Cells.Item(193, 2).Value = 101
$ws.Cells.Item(193, 3).Value = 4.26000022888184
$ws.Cells.Item(193, 4).Value = 4.11999988555908
$ws.Cells.Item(193, 5).Value = 4.25
$ws.Cells.Item(193, 6).Value = 4.11999988555908
$ws.Cells.Item(193, 7).NumberFormat = "@"
$ws.Cells.Item(193, 7).Value = "4.11999988555908"
$ws.Cells.Item(193, 7).Style = "Normal"
$ws.Cells.Item(193, 8).Value = "ESF.MI"

# Column A uses the custom date/time number format applied throughout the
# sheet (style index 1 in the original file) - copy that style down from
# the row above so the two new rows match the rest of the column.
$ws.Range("A191").Copy()
$ws.Range("A192:A193").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
